$d = $word.ActiveDocument

# Locate the paragraph that ends with the "LOM3246..." text, then delete
# the three paragraphs that follow it: the blank paragraph, the
# "Ver no Jupiter..." paragraph, and the "(c) 2020 ..." paragraph.
# The blank paragraph and page-break paragraph that come after those
# three must be left untouched.

$count = $d.Paragraphs.Count
$startIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOM3246") {
        $startIndex = $i
        break
    }
}

if ($startIndex -gt 0) {
    # Build a range spanning the three paragraphs right after the LOM3246 one
    # (indexes startIndex+1, startIndex+2, startIndex+3) and delete it.
    $first = $d.Paragraphs.Item($startIndex + 1)
    $last = $d.Paragraphs.Item($startIndex + 3)

    $deleteRange = $d.Range($first.Range.Start, $last.Range.End)
    $deleteRange.Delete()
}
